# categories.xlsx export refresh: rename status columns, refresh the
# isActive flag into a createdAt/updatedAt/status layout, and append the
# newer category rows (withdrawalMethode + feedback changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row: isActive/createdAt/updatedAt -> createdAt/updatedAt/status
# ---------------------------------------------------------------------
$ws.Range("E1").Value = "createdAt"
$ws.Range("F1").Value = "updatedAt"
$ws.Range("G1").Value = "status"

# ---------------------------------------------------------------------
# A date-formatted "template" cell (existing F2) that we clone the date
# number-format from, so new/changed date cells reuse the same style
# (s="1") instead of minting a duplicate numFmt.
# ---------------------------------------------------------------------
$dateTemplate = $ws.Range("F2")

function Set-DateCell($cell, $serial) {
    $cell.Value = $serial
    $dateTemplate.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

function Set-StatusCell($cell, $text) {
    # These cells previously held the "updatedAt" date (date number format);
    # drop that formatting so the new status text isn't tagged with it.
    $cell.ClearFormats()
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# Row 2 (Fashion): E2 old "active" flag -> createdAt date; F2 createdAt
# date -> updatedAt date (new value); G2 updatedAt date -> status text.
# ---------------------------------------------------------------------
Set-DateCell $ws.Range("E2") 45581.75886135417
$ws.Range("F2").Value = 45584.72708806713
Set-StatusCell $ws.Range("G2") "Active"

# ---------------------------------------------------------------------
# Row 3 (Home Appliences)
# ---------------------------------------------------------------------
Set-DateCell $ws.Range("E3") 45581.76085006945
$ws.Range("F3").Value = 45584.508307986114
Set-StatusCell $ws.Range("G3") "Active"

# ---------------------------------------------------------------------
# Row 4: category document was replaced (Cloths -> Furniture, new id/logo)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "6719f2bb49e1f74e65127577"
$ws.Range("B4").Value = "Furniture"
$ws.Range("D4").Value = "https://res.cloudinary.com/dgexhjryd/image/upload/v1729753811/Lailoji/logo-1729753810993.jpg"
Set-DateCell $ws.Range("E4") 45589.52774784722
$ws.Range("F4").Value = 45589.52831390046
Set-StatusCell $ws.Range("G4") "Inactive"

# ---------------------------------------------------------------------
# New rows 5-11: categories added since the last export
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "671b2a00e944f0b5198daf9a"
$ws.Range("B5").Value = "Handbags"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "https://res.cloudinary.com/dgexhjryd/image/upload/v1729082555/Lailoji/logo-1729082551627.jpg"
Set-DateCell $ws.Range("E5") 45590.45002712963
Set-DateCell $ws.Range("F5") 45590.45002712963
Set-StatusCell $ws.Range("G5") "Inactive"

$ws.Range("A6").Value = "671b2a00e944f0b5198daf9b"
$ws.Range("B6").Value = "Shoes"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "https://res.cloudinary.com/dgexhjryd/image/upload/v1729082555/Lailoji/logo-1729082551627.jpg"
Set-DateCell $ws.Range("E6") 45590.45002712963
Set-DateCell $ws.Range("F6") 45590.45002712963
Set-StatusCell $ws.Range("G6") "Inactive"

$ws.Range("A7").Value = "6729ac9fbe93e05ec4836723"
$ws.Range("B7").Value = "Garden"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "https://res.cloudinary.com/dgexhjryd/image/upload/v1730784463/Lailoji/logo-1730784463444.jpg"
Set-DateCell $ws.Range("E7") 45601.45631443287
Set-DateCell $ws.Range("F7") 45601.45687869213
Set-StatusCell $ws.Range("G7") "Inactive"

$ws.Range("A8").Value = "6729b3472ea45155a9a73ffb"
$ws.Range("B8").Value = "Electronics"
$ws.Range("C8").Value = "high"
$ws.Range("D8").Value = "https://res.cloudinary.com/dgexhjryd/image/upload/v1730786356/Lailoji/logo-1730786353820.jpg"
Set-DateCell $ws.Range("E8") 45601.47604179398
Set-DateCell $ws.Range("F8") 45601.47881064815
Set-StatusCell $ws.Range("G8") "Inactive"

$ws.Range("A9").Value = "6729bbcf8ac6688c38353ce9"
$ws.Range("B9").Value = "Medicine"
$ws.Range("C9").Value = "low"
$ws.Range("D9").Value = "https://res.cloudinary.com/dgexhjryd/image/upload/v1730788301/Lailoji/logo-1730788300510.jpg"
Set-DateCell $ws.Range("E9") 45601.501319375
Set-DateCell $ws.Range("F9") 45601.501319375
Set-StatusCell $ws.Range("G9") "Inactive"

$ws.Range("A10").Value = "6729bc308ac6688c38353cec"
$ws.Range("B10").Value = "gif"
$ws.Range("C10").Value = "low"
$ws.Range("D10").Value = "https://res.cloudinary.com/dgexhjryd/image/upload/v1730788398/Lailoji/logo-1730788397696.gif"
Set-DateCell $ws.Range("E10") 45601.50244056713
Set-DateCell $ws.Range("F10") 45601.50244056713
Set-StatusCell $ws.Range("G10") "Inactive"

$ws.Range("A11").Value = "6729be8d8ac6688c38353cfd"
$ws.Range("B11").Value = "Fruits"
$ws.Range("C11").Value = "low"
$ws.Range("D11").Value = "https://res.cloudinary.com/dgexhjryd/image/upload/v1730789069/Lailoji/logo-1730789069069.jpg"
Set-DateCell $ws.Range("E11") 45601.50943373843
Set-DateCell $ws.Range("F11") 45601.5102078125
Set-StatusCell $ws.Range("G11") "Inactive"

# ---------------------------------------------------------------------
# Keep the "number stored as text" ignore-range in sync with the grown
# data extent (A1:G4 -> A1:G11), mirroring what Excel's error-checking
# UI does when you dismiss the warning over the full used range.
# ---------------------------------------------------------------------
try {
    $ws.Range("A1:G11").Errors.Item(9).Ignore = $true
} catch {
}
